$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update per-row Price (D, stored as text) and Volume(1h) (E) values.
# Force NumberFormat to Text on Price cells first so Excel does not
# reinterpret numeric-looking strings (e.g. "1.00", "214.00") as numbers
# and strip the trailing zeros / formatting, matching the original
# inline-string (text) cell content exactly.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.866.81'
$ws.Range("E2").Value = '  +0.77%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.626.78'
$ws.Range("E3").Value = '  +1.50%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.00'
$ws.Range("E5").Value = '  +0.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.520'
$ws.Range("E6").Value = '  +0.60%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '29.33'
$ws.Range("E8").Value = '  +8.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.258'
$ws.Range("E9").Value = '  +2.81%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0611'
$ws.Range("E10").Value = '  +1.79%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0914'
$ws.Range("E11").Value = '  +0.84%  '
$ws.Range("E12").Value = '  +1.49%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.626.27'
$ws.Range("E13").Value = '  +1.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.572'
$ws.Range("E14").Value = '  +6.77%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.89'
$ws.Range("E15").Value = '  +4.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.940.76'
$ws.Range("E16").Value = '  +0.99%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '8.88'
$ws.Range("E17").Value = '  +17.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.65'
$ws.Range("E18").Value = '  +2.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.89'
$ws.Range("E19").Value = '  +0.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0705'
$ws.Range("E20").Value = '  +2.02%  '
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.13'
$ws.Range("E22").Value = '  +3.47%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.62'
$ws.Range("E23").Value = '  +4.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.12'
$ws.Range("E24").Value = '  +2.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.54'
$ws.Range("E25").Value = '  +1.92%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.69'
$ws.Range("E26").Value = '  +2.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.110'
$ws.Range("E27").Value = '  +2.70%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.61'
$ws.Range("E28").Value = '  +3.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0489'
$ws.Range("E30").Value = '  +3.00%  '
$ws.Range("E31").Value = '  +5.39%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.34'
$ws.Range("E32").Value = '  +3.81%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.21'
$ws.Range("E33").Value = '  +2.69%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.423.04'
$ws.Range("E34").Value = '  -0.55%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.63'
$ws.Range("E35").Value = '  +6.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.03'
$ws.Range("E36").Value = '  +0.26%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.87'
$ws.Range("E37").Value = '  +1.87%  '
$ws.Range("E38").Value = '  -0.35%  '
$ws.Range("E39").Value = '  +3.32%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.557'
$ws.Range("E40").Value = '  +3.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0498'
$ws.Range("E43").Value = '  +2.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '54.23'
$ws.Range("E44").Value = '  -1.35%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '69.46'
$ws.Range("E45").Value = '  +5.55%  '
$ws.Range("E46").Value = '  +6.73%  '
$ws.Range("E47").Value = '  +0.09%  '
$ws.Range("E48").Value = '  +2.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.768.10'
$ws.Range("E49").Value = '  +1.35%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '88.69'
$ws.Range("E50").Value = '  +2.62%  '
$ws.Range("E51").Value = '  +6.23%  '

# Rows 41 and 42 swap coin identity (ARBITRUM <-> RenderToken) with updated data
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.830'
$ws.Range("E41").Value = '  +3.79%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.98'
$ws.Range("E42").Value = '  +0.55%  '
